# [Fonds de solidarite] Add 2020-12-14 data
# Updates nombre_aides (C), nombre_entreprises (D) and montant_total (E)
# for the rows whose underlying source figures changed with the new data.
#
# The sheet stores every value as text (inlineStr / shared string), even
# though many of them look numeric. Writing plain numeric-looking strings
# through Range.Value would make Excel auto-convert the cell to a real
# Number, which would both change the cell's type and risk losing exact
# textual formatting (e.g. trailing zeros like "530310.00"). To avoid
# that, each touched cell is temporarily forced to Text format before the
# write, then restored to the default "Normal" style so no stray
# NumberFormat/style is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Value
    )
    $rng = $ws.Range($Address)
    $rng.NumberFormat = "@"
    $rng.Value = $Value
    $rng.Style = "Normal"
}

# Row 36
Set-TextValue "C36" "572"
Set-TextValue "E36" "3831528.78"

# Row 49
Set-TextValue "C49" "833"
Set-TextValue "E49" "2522695.54"

# Row 52
Set-TextValue "C52" "812"
Set-TextValue "E52" "5537284.76"

# Row 63
Set-TextValue "C63" "156"
Set-TextValue "D63" "156"
Set-TextValue "E63" "374491.16"

# Row 67
Set-TextValue "C67" "193"
Set-TextValue "D67" "189"
Set-TextValue "E67" "530310.00"

# Row 68
Set-TextValue "C68" "260"
Set-TextValue "D68" "257"
Set-TextValue "E68" "873089.81"

# Row 69
Set-TextValue "C69" "197"
Set-TextValue "D69" "192"
Set-TextValue "E69" "563604.35"
